$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23/24 swap: Chainlink <-> BinanceUSD reordering
$ws.Range("B23").Value = 'Chainlink'
$ws.Range("B24").Value = 'BinanceUSD'
$ws.Range("C23").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("C24").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'

# Price (column D) values: many look numeric, so Excel would silently coerce a plain
# .Value assignment into a Number (losing formatting like trailing zeros / long decimals).
# Instead: write each as a text-formula ("=""..."""), then Copy + PasteSpecial(xlPasteValues)
# over the same block. That freezes the formula result as a literal, keeping it a String
# without touching NumberFormat/style on the cell.
$ws.Range("D2").Formula = '="26.036.66"'
$ws.Range("D3").Formula = '="1.665.47"'
$ws.Range("D4").Formula = '="1.006"'
$ws.Range("D5").Formula = '="217.83"'
$ws.Range("D6").Formula = '="0.5038"'
$ws.Range("D7").Formula = '="1.007"'
$ws.Range("D8").Formula = '="0.2633"'
$ws.Range("D9").Formula = '="0.06318"'
$ws.Range("D10").Formula = '="21.42"'
$ws.Range("D11").Formula = '="0.07360"'
$ws.Range("D12").Formula = '="1.671.05"'
$ws.Range("D13").Formula = '="4.543"'
$ws.Range("D14").Formula = '="0.5736"'
$ws.Range("D15").Formula = '="1.891.39"'
$ws.Range("D16").Formula = '="0.000008428"'
$ws.Range("D17").Formula = '="64.52"'
$ws.Range("D18").Formula = '="26.090.31"'
$ws.Range("D19").Formula = '="4.932"'
$rng = $ws.Range("D2:D19")
$rng.Copy()
$rng.PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("D21").Formula = '="10.75"'
$ws.Range("D22").Formula = '="186.52"'
$ws.Range("D23").Formula = '="6.164"'
$ws.Range("D24").Formula = '="1.007"'
$ws.Range("D25").Formula = '="142.62"'
$ws.Range("D26").Formula = '="7.653"'
$ws.Range("D27").Formula = '="0.1168"'
$ws.Range("D28").Formula = '="15.73"'
$ws.Range("D29").Formula = '="1.296"'
$ws.Range("D30").Formula = '="0.05783"'
$ws.Range("D31").Formula = '="1.322"'
$ws.Range("D32").Formula = '="3.490"'
$ws.Range("D33").Formula = '="3.498"'
$ws.Range("D34").Formula = '="1.646"'
$ws.Range("D35").Formula = '="1.003"'
$ws.Range("D36").Formula = '="0.5968"'
$ws.Range("D37").Formula = '="2.364"'
$ws.Range("D38").Formula = '="2.638"'
$ws.Range("D39").Formula = '="0.01595"'
$ws.Range("D40").Formula = '="1.081.54"'
$ws.Range("D41").Formula = '="5.949"'
$ws.Range("D42").Formula = '="0.8587"'
$rng = $ws.Range("D21:D42")
$rng.Copy()
$rng.PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("D44").Formula = '="99.46"'
$ws.Range("D45").Formula = '="1.810.56"'
$ws.Range("D46").Formula = '="0.00000000113"'
$ws.Range("D47").Formula = '="55.80"'
$ws.Range("D48").Formula = '="1.007"'
$ws.Range("D49").Formula = '="8.107"'
$ws.Range("D50").Formula = '="0.4298"'
$ws.Range("D51").Formula = '="0.05168"'
$rng = $ws.Range("D44:D51")
$rng.Copy()
$rng.PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Volume(1h) (column E) values: always wrapped in spaces + "%", never auto-converted.
$ws.Range("E2").Value = '  -7.17%  '
$ws.Range("E3").Value = '  -4.59%  '
$ws.Range("E4").Value = '  +0.51%  '
$ws.Range("E5").Value = '  -3.91%  '
$ws.Range("E6").Value = '  -12.99%  '
$ws.Range("E7").Value = '  +0.45%  '
$ws.Range("E8").Value = '  -3.01%  '
$ws.Range("E9").Value = '  -4.33%  '
$ws.Range("E10").Value = '  -7.44%  '
$ws.Range("E12").Value = '  -4.34%  '
$ws.Range("E13").Value = '  -3.86%  '
$ws.Range("E14").Value = '  -5.20%  '
$ws.Range("E15").Value = '  -4.59%  '
$ws.Range("E16").Value = '  -2.79%  '
$ws.Range("E17").Value = '  -13.14%  '
$ws.Range("E18").Value = '  -6.94%  '
$ws.Range("E19").Value = '  -7.59%  '
$ws.Range("E20").Value = '  +0.55%  '
$ws.Range("E21").Value = '  -4.71%  '
$ws.Range("E22").Value = '  -9.04%  '
$ws.Range("E23").Value = '  -7.09%  '
$ws.Range("E24").Value = '  +0.39%  '
$ws.Range("E25").Value = '  -4.90%  '
$ws.Range("E26").Value = '  -5.20%  '
$ws.Range("E27").Value = '  -5.34%  '
$ws.Range("E28").Value = '  -2.45%  '
$ws.Range("E29").Value = '  -6.65%  '
$ws.Range("E30").Value = '  -7.30%  '
$ws.Range("E31").Value = '  -4.95%  '
$ws.Range("E32").Value = '  -6.72%  '
$ws.Range("E33").Value = '  -5.80%  '
$ws.Range("E34").Value = '  -1.91%  '
$ws.Range("E35").Value = '  -3.27%  '
$ws.Range("E36").Value = '  -6.34%  '
$ws.Range("E37").Value = '  -3.65%  '
$ws.Range("E38").Value = '  -3.17%  '
$ws.Range("E39").Value = '  -4.70%  '
$ws.Range("E40").Value = '  -4.11%  '
$ws.Range("E41").Value = '  -4.02%  '
$ws.Range("E42").Value = '  -1.76%  '
$ws.Range("E43").Value = '  +0.40%  '
$ws.Range("E44").Value = '  -0.17%  '
$ws.Range("E45").Value = '  -4.47%  '
$ws.Range("E46").Value = '  +4.80%  '
$ws.Range("E47").Value = '  -6.15%  '
$ws.Range("E48").Value = '  +0.70%  '
$ws.Range("E49").Value = '  -1.99%  '
$ws.Range("E50").Value = '  -2.78%  '
$ws.Range("E51").Value = '  -3.98%  '
